$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 300, shifting existing rows 300-373 down to 301-374.
$ws.Rows.Item(300).Insert()

# Populate the newly inserted row 300 with the new data record.
$ws.Cells.Item(300, 1).Value = 7
$ws.Cells.Item(300, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(300, 3).Value = "Ñuble"
$ws.Cells.Item(300, 4).Value = 45173
$ws.Cells.Item(300, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(300, 5).Value = 16
$ws.Cells.Item(300, 6).Value = 100112045
$ws.Cells.Item(300, 7).Value = "Zapallo"
$ws.Cells.Item(300, 8).Value = "Paine"
$ws.Cells.Item(300, 9).Value = "1a (guarda)"
$ws.Cells.Item(300, 10).Value = 250
$ws.Cells.Item(300, 11).Value = 350
$ws.Cells.Item(300, 12).Value = 350
$ws.Cells.Item(300, 13).Value = 350
$ws.Cells.Item(300, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(300, 15).Value = "Región del Maule"
$ws.Cells.Item(300, 16).Value = 350
$ws.Cells.Item(300, 17).Value = 1
$ws.Cells.Item(300, 18).Value = "Hortaliza"
